$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Información Modelo" sheet: update objective value, gap time and
#    constraint count.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Información Modelo")
$wsInfo.Range("A2").Value = 640108574274.0112
$wsInfo.Range("C2").Value = 1.584
$wsInfo.Range("E2").Value = 24926

# ---------------------------------------------------------------------------
# 2) "Procesos Activados" sheet: replace the 3 data rows with 19 new rows,
#    column A is always 1, column B steps by 20 (0, 20, 40, ..., 360).
# ---------------------------------------------------------------------------
$wsProcAct = $wb.Worksheets.Item("Procesos Activados")
for ($i = 0; $i -lt 19; $i++) {
    $r = $i + 2
    $wsProcAct.Cells.Item($r, 1).Value = 1
    $wsProcAct.Cells.Item($r, 2).Value = $i * 20
}

# ---------------------------------------------------------------------------
# 3) "Procesos en Operación" sheet: column A changes from 4 to 1 on every
#    data row (rows 2..366); column B (Tiempo) is unchanged.
# ---------------------------------------------------------------------------
$wsProcOp = $wb.Worksheets.Item("Procesos en Operación")
for ($r = 2; $r -le 366; $r++) {
    $wsProcOp.Cells.Item($r, 1).Value = 1
}

# ---------------------------------------------------------------------------
# 4) "Total Contaminantes Z" sheet: update totals per contaminant.
# ---------------------------------------------------------------------------
$wsTotalZ = $wb.Worksheets.Item("Total Contaminantes Z")
$wsTotalZ.Range("B2").Value = 449208244800.0004
$wsTotalZ.Range("B3").Value = 13481640000.00001
$wsTotalZ.Range("B4").Value = 87091394399.99998
$wsTotalZ.Range("B5").Value = 307074.010608
$wsTotalZ.Range("B6").Value = 90326988000.00008

# ---------------------------------------------------------------------------
# 5) "Concentraciones" sheet: update output concentration per contaminant.
# ---------------------------------------------------------------------------
$wsConc = $wb.Worksheets.Item("Concentraciones")
$wsConc.Range("B2").Value = 16.66
$wsConc.Range("B3").Value = 0.5
$wsConc.Range("B4").Value = 3.23
$wsConc.Range("B5").Value = 0.000011
$wsConc.Range("B6").Value = 3.35

# ---------------------------------------------------------------------------
# 6) New "Costo Total" sheet, appended as the last tab, carrying the single
#    total-cost result value. Its header cell reuses the same bold/bordered
#    style used by the other sheets' header rows.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCosto = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCosto.Name = "Costo Total"

$wsInfo.Range("A1").Copy($wsCosto.Range("A1"))
$wsCosto.Range("A1").Value = "Costo Total"
$wsCosto.Range("A2").Value = 1136.135
